# The "Description" column (column B) is being removed from the
# defect-upload template. Deleting it shifts the "Fields marked *...
# are required." note (previously column C) into column B, and the
# previously-empty column D into column C.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = $ws.Range("B1:B1048576")
$colB.Select()
$colB.EntireColumn.Delete()
